# Update the "Förändrad" (changed/updated) date column C for every data
# row (rows 2-203) from 2023-09-15 (45184) to 2023-09-17 (45186).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C203").Value = 45186

# Add a friendly display-text second argument to every HYPERLINK() formula
# in the sheet (columns S, T, U, V, W, X, Y) so the link shows the
# "Beteckning" (case id, column A) instead of the raw URL. Only touch
# cells that actually contain a HYPERLINK(...) formula and don't already
# carry a second argument.
$cols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = 2; $row -le 203; $row++) {
    $caseId = $ws.Range("A" + $row).Value2
    if (-not $caseId) { continue }

    foreach ($col in $cols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula

        if ($formula -and $formula.Length -ge 11 -and $formula.Substring(0, 11) -eq "=HYPERLINK(" -and -not $formula.Contains(",")) {
            $updated = $formula.Substring(0, $formula.Length - 1) + ', "' + $caseId + '")'
            $cell.Formula = $updated
        }
    }
}
